# Update Leve profit-tracking figures across all profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2850.1
$ws.Range("I70").Value = 2850.1
$ws.Range("K70").Value = 8550.299999999999
$ws.Range("M70").Value = -8280.299999999999
$ws.Range("H73").Value = 2850.1
$ws.Range("I73").Value = 2850.1
$ws.Range("K73").Value = 8550.299999999999
$ws.Range("M73").Value = -7614.299999999999
$ws.Range("H98").Value = 11913.272
$ws.Range("I98").Value = 5550.25
$ws.Range("K98").Value = 5550.25
$ws.Range("M98").Value = -4052.25
$ws.Range("H107").Value = 1790.1666
$ws.Range("I107").Value = 1079.1666
$ws.Range("J107").Value = 2501.1667
$ws.Range("K107").Value = 1079.1666
$ws.Range("L107").Value = 2501.1667
$ws.Range("M107").Value = 840.8334
$ws.Range("N107").Value = -6341.1667
$ws.Range("H113").Value = 47500
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -6746
$ws.Range("H116").Value = 2750
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("M116").Value = 1442
$ws.Range("H122").Value = 11913.272
$ws.Range("I122").Value = 5550.25
$ws.Range("K122").Value = 16650.75
$ws.Range("M122").Value = -14200.75
$ws.Range("H132").Value = 2900.1428
$ws.Range("I132").Value = 860.6
$ws.Range("K132").Value = 2581.8
$ws.Range("M132").Value = -51.80000000000018
$ws.Range("H137").Value = 1074.5
$ws.Range("I137").Value = 1074.5
$ws.Range("K137").Value = 3223.5
$ws.Range("M137").Value = -673.5
$ws.Range("H141").Value = 3760.75
$ws.Range("I141").Value = 3596.2
$ws.Range("K141").Value = 10788.6
$ws.Range("M141").Value = -5608.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 708.75
$ws.Range("I32").Value = 708.75
$ws.Range("K32").Value = 708.75
$ws.Range("M32").Value = -421.75
$ws.Range("H61").Value = 1774.375
$ws.Range("I61").Value = 1742.1428
$ws.Range("K61").Value = 1742.1428
$ws.Range("M61").Value = -1530.1428
$ws.Range("H132").Value = 1708.3334
$ws.Range("I132").Value = 1553
$ws.Range("K132").Value = 4659
$ws.Range("M132").Value = -2129
$ws.Range("H136").Value = 1774.375
$ws.Range("I136").Value = 1742.1428
$ws.Range("K136").Value = 5226.428400000001
$ws.Range("M136").Value = -2676.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 76666.664
$ws.Range("J13").Value = 76666.664
$ws.Range("L13").Value = 76666.664
$ws.Range("N13").Value = -77002.664
$ws.Range("H134").Value = 2228
$ws.Range("I134").Value = 1970.6666
$ws.Range("K134").Value = 5911.9998
$ws.Range("M134").Value = -3376.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16798.7
$ws.Range("I31").Value = 7166.5
$ws.Range("J31").Value = 31247
$ws.Range("K31").Value = 7166.5
$ws.Range("L31").Value = 31247
$ws.Range("M31").Value = -6871.5
$ws.Range("N31").Value = -31837
$ws.Range("H34").Value = 16798.7
$ws.Range("I34").Value = 7166.5
$ws.Range("J34").Value = 31247
$ws.Range("K34").Value = 7166.5
$ws.Range("L34").Value = 31247
$ws.Range("M34").Value = -6964.5
$ws.Range("N34").Value = -31651
$ws.Range("H98").Value = 28625
$ws.Range("I98").Value = 25000
$ws.Range("J98").Value = 32250
$ws.Range("K98").Value = 25000
$ws.Range("L98").Value = 32250
$ws.Range("M98").Value = -22754
$ws.Range("N98").Value = -36742
$ws.Range("H134").Value = 1278.7778
$ws.Range("I134").Value = 1244.1428
$ws.Range("K134").Value = 3732.4284
$ws.Range("M134").Value = -1197.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 558.4
$ws.Range("J5").Value = 705
$ws.Range("L5").Value = 2115
$ws.Range("N5").Value = -2339
$ws.Range("H44").Value = 1035.7142
$ws.Range("I44").Value = 450
$ws.Range("J44").Value = 2500
$ws.Range("K44").Value = 1350
$ws.Range("L44").Value = 7500
$ws.Range("M44").Value = -952
$ws.Range("N44").Value = -8296
$ws.Range("H113").Value = 196.11111
$ws.Range("I113").Value = 249.66667
$ws.Range("J113").Value = 169.33333
$ws.Range("K113").Value = 749.00001
$ws.Range("L113").Value = 507.99999
$ws.Range("M113").Value = 1420.99999
$ws.Range("N113").Value = -4847.99999
$ws.Range("H135").Value = 558.4
$ws.Range("J135").Value = 705
$ws.Range("L135").Value = 6345
$ws.Range("N135").Value = -11415

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3576.923
$ws.Range("I36").Value = 2700
$ws.Range("K36").Value = 2700
$ws.Range("M36").Value = -2215
$ws.Range("H43").Value = 1820
$ws.Range("I43").Value = 1500
$ws.Range("K43").Value = 1500
$ws.Range("M43").Value = -1349
$ws.Range("H122").Value = 4631.222
$ws.Range("I122").Value = 3375
$ws.Range("J122").Value = 7143.6665
$ws.Range("K122").Value = 10125
$ws.Range("L122").Value = 21430.9995
$ws.Range("M122").Value = -7675
$ws.Range("N122").Value = -26330.9995
$ws.Range("H123").Value = 60333.5
$ws.Range("J123").Value = 60333.5
$ws.Range("L123").Value = 60333.5
$ws.Range("N123").Value = -65233.5
$ws.Range("H126").Value = 11998.363
$ws.Range("I126").Value = 11198.2
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 33594.60000000001
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -31124.60000000001
$ws.Range("N126").Value = -64940
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2200
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6600
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -4070
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 299999
$ws.Range("J5").Value = 299999
$ws.Range("L5").Value = 299999
$ws.Range("N5").Value = -300225
$ws.Range("H40").Value = 25249.363
$ws.Range("I40").Value = 20968
$ws.Range("J40").Value = 36666.332
$ws.Range("K40").Value = 20968
$ws.Range("L40").Value = 36666.332
$ws.Range("M40").Value = -20832
$ws.Range("N40").Value = -36938.332
$ws.Range("H46").Value = 287883
$ws.Range("I46").Value = 668000
$ws.Range("J46").Value = 2795.25
$ws.Range("K46").Value = 668000
$ws.Range("L46").Value = 2795.25
$ws.Range("M46").Value = -667812
$ws.Range("N46").Value = -3171.25
$ws.Range("H55").Value = 279.625
$ws.Range("J55").Value = 240.25
$ws.Range("L55").Value = 240.25
$ws.Range("N55").Value = -586.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8651.77
$ws.Range("I136").Value = 7580.3335
$ws.Range("J136").Value = 9570.143
$ws.Range("K136").Value = 22741.0005
$ws.Range("L136").Value = 28710.429
$ws.Range("M136").Value = -20191.0005
$ws.Range("N136").Value = -33810.429
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = $null
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null
